# Update homepage publications: add the newly-accepted Bifrost (VLDB 2026)
# and PMark (ICLR 2026, moved from Preprints) papers to the Publications
# sheet, and drop the now-published PMark entry from the Preprints sheet.

$wb = $excel.ActiveWorkbook

# --- Publications sheet: insert two new rows right after the header ---
$ws1 = $wb.Worksheets.Item("Publications")

$ws1.Rows("2:3").Insert()
$ws1.Rows("2:3").ClearFormats()

$ws1.Range("A2").Value = 2026
$ws1.Range("B2").Value = "VLDB"
$ws1.Range("C2").Value = "Bifrost: A Much Simpler Secure Two-Party Data Join Protocol for Secure Data Analytics"
$ws1.Range("D2").Value = "https://wuwuz.github.io"
$ws1.Range("E2").Value = "Shuyu Chen, Mingxun Zhou, Haoyu Niu, Guopeng Lin, Weili Han"
$ws1.Range("F2:G2").ClearContents()

$ws1.Range("A3").Value = 2026
$ws1.Range("B3").Value = "ICLR"
$ws1.Range("C3").Value = "PMark: Towards Robust and Distortion-free Semantic-level Watermarking with Channel Constraints"
$ws1.Range("D3").Value = "https://arxiv.org/abs/2509.21057"
$ws1.Range("E3").Value = "Jiahao Huo, Shuliang Liu, Bin Wang, Junyan Zhang, Yibo Yan, Aiwei Liu, Xuming Hu, Mingxun Zhou"
$ws1.Range("F3:G3").ClearContents()

# --- Preprints sheet: the PMark entry is now published, remove it ---
$ws2 = $wb.Worksheets.Item("Preprints")
$ws2.Rows("2:2").Delete()

# --- View state: Publications selection moves, Preprints becomes the active tab ---
$ws1.Activate()
$ws1.Range("C30").Select()

$ws2.Activate()
$ws2.Range("A11").Select()
